$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) / "Valor Mora" (column F) rows 16-21 are being
# re-sorted so periods run in ascending order (1801 -> 1806) instead of the
# previous descending order (1806 -> 1801). Only the displayed values move;
# the per-row formatting (borders etc.) stays where it is, so we set cell
# values directly instead of doing a range sort (which would also swap the
# formatting between rows).

$periods = @("1801", "1802", "1803", "1804", "1805", "1806")
$valores = @(16666, 31249, 31249, 31249, 31249, 31249)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
